$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.444.33'
$ws.Range('E2').Value = '  +4.28%  '
$ws.Range('D3').Value = '3.003.63'
$ws.Range('E3').Value = '  +2.76%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''562.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('D6').Value = '''138.35'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.60%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.523'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.60%  '
$ws.Range('D9').Value = '2.987.57'
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('D11').Value = '''5.16'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.37%  '
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('E13').Value = '  +5.45%  '
$ws.Range('D14').Value = '''33.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.28%  '
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').Value = '3.499.32'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('E17').Value = '  +7.76%  '
$ws.Range('D18').Value = '2.998.73'
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('D19').Value = '59.405.40'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('D20').Value = '''430.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('D21').Value = '''13.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.78%  '
$ws.Range('E22').Value = '  +6.15%  '
$ws.Range('D23').Value = '''7.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.96%  '
$ws.Range('D24').Value = '''13.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.59%  '
$ws.Range('D25').Value = '''80.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.56%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').Value = '''2.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.51%  '
$ws.Range('E29').Value = '  +2.98%  '
$ws.Range('D30').Value = '''7.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.36%  '
$ws.Range('D31').Value = '''25.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.39%  '
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = '''0.0993'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').Value = '''0.994'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.89%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '''5.92'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.85%  '
$ws.Range('D36').Value = '0.0₃0763'
$ws.Range('E36').Value = '  +14.22%  '
$ws.Range('E37').Value = '  +0.91%  '
$ws.Range('D38').Value = '''49.07'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('D40').Value = '''2.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.88%  '
$ws.Range('D41').Value = '''400.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.34%  '
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('D43').Value = '2.753.09'
$ws.Range('E43').Value = '  +4.68%  '
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').Value = '''0.252'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.07%  '
$ws.Range('D46').Value = '''35.79'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +27.93%  '
$ws.Range('D47').Value = '''0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').Value = '''122.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('E51').Value = '  +1.40%  '
